$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the "Data" column text: "2024-23-06" -> " 2024-06-23" (leading space preserved)
#    for every data row (D2:D272), all of which shared the same string.
$ws.Range("D2:D272").Value = " 2024-06-23"

# 2. Format that same column as a date (built-in numFmtId 14) - apply to one cell
#    then copy/paste the formatting onto the rest so every cell shares a single
#    style record instead of one per cell.
$ws.Range("D2").NumberFormat = "mm-dd-yy"
$null = $ws.Range("D2").Copy()
$null = $ws.Range("D3:D272").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Update the active selection shown when the sheet is opened.
$null = $ws.Range("E10").Select()

# 4. Configure the print page setup (A4, portrait).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
